# Update "想去人数" (number of people interested) counts on the
# "展览" and "全部类型" sheets, as produced by the site's scheduled
# data refresh (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览": rows 2-5 hold F column counts.
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 290
$wsExpo.Range("F3").Value = 1109
$wsExpo.Range("F4").Value = 2535
$wsExpo.Range("F5").Value = 217

# Sheet "全部类型": same underlying events, different rows.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 290
$wsAll.Range("F5").Value = 1109
$wsAll.Range("F6").Value = 2535
$wsAll.Range("F8").Value = 217
